$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Upcoming Events" heading (E2) to "Upcoming Events!"
$ws.Range("E2").Value = "Upcoming Events!"

# Update the "Thursdays!" heading (E3) to "Every Thursday!"
$ws.Range("E3").Value = "Every Thursday!"

# Update the selected cell in the sheet view to match the saved state
$ws.Range("D13").Select()
